# Added filter per company and excel generation
#
# 1. Rename "First" -> "Acosta, Butler and Perez"
# 2. Remove the now-unused "New" sheet
# 3. Rebuild the main data table on the first sheet with the new,
#    wider schema (14 columns) and new sample rows
# 4. Re-activate the "Evaluation Warning" sheet (keeps it the selected tab)

$wb = $excel.ActiveWorkbook

# --- 1. rename first sheet -------------------------------------------------
$wb.Worksheets.Item("First").Name = "Acosta, Butler and Perez"

# --- 2. drop the "New" sheet ------------------------------------------------
$wb.Worksheets.Item("New").Delete()

# --- 3. rebuild the data table ----------------------------------------------
$ws = $wb.Worksheets.Item("Acosta, Butler and Perez")
$ws.Cells.Clear()

$headers = @("Name","College ID","Contact number","Email","Date of Birth","Gender","Program","Specialization","10th Marks","12th Marks","CGPA","Backlogs","Red Flags","Category")

$row2 = @("Inna","2023PCP5321","'8050106439","niharkajla123@gmail.com","'2024-01-16","Female","PG","CE","'90","'90","'9","'0","'0","SC")

$row3 = @("Suhana Sharma","2023PCP5305","'8050106439","niharamazon5005@gmail.com","'2024-03-07","Female","PG","CSE","'90","'90","'9","'0","'0","General")

for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}
for ($c = 0; $c -lt $row2.Length; $c++) {
    $ws.Cells.Item(2, $c + 1).Value = $row2[$c]
}
for ($c = 0; $c -lt $row3.Length; $c++) {
    $ws.Cells.Item(3, $c + 1).Value = $row3[$c]
}

# header row is bold, like the original header
$ws.Range("A1:N1").Font.Bold = $true

# drop the now-unused trailing rows that used to hold rows 4-7 of data
$ws.Range("A4:N7").EntireRow.Delete()

# --- 4. keep "Evaluation Warning" as the active / selected sheet -----------
$wb.Worksheets.Item("Evaluation Warning").Activate()
